# Append the 2024-01-04 09:30 resale-number row (row 15) to the
# CityResaleNum sheet, matching the columns already present in row 14.
#
# Columns A-D are text (Date/Time/Weekday/Week), columns E-T are numbers.
# "2024-01-04" and "00" both look numeric/date-like to Excel's type
# inference, so those two cells are pre-formatted as Text to keep them
# as literal strings instead of being auto-converted to a date serial /
# stripped-leading-zero number (matches the existing rows' values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

$ws.Range("A$row").NumberFormat = "@"
$ws.Range("D$row").NumberFormat = "@"

$ws.Range("A$row").Value = "2024-01-04"
$ws.Range("B$row").Value = "09:30:22"
$ws.Range("C$row").Value = "Thursday"
$ws.Range("D$row").Value = "00"

$ws.Range("E$row").Value = 140160
$ws.Range("F$row").Value = 142893
$ws.Range("G$row").Value = 171703
$ws.Range("H$row").Value = 146644
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 117287
$ws.Range("K$row").Value = 223915
$ws.Range("L$row").Value = 247889
$ws.Range("M$row").Value = 184235
$ws.Range("N$row").Value = 109943
$ws.Range("O$row").Value = 40216
$ws.Range("P$row").Value = 30833
$ws.Range("Q$row").Value = 72211
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 41035
$ws.Range("T$row").Value = -1
